$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column R: year header "2021" (copy formatting from Q3) ---
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)
$ws.Range("R3").Value = 2021

# --- Add new column R: data value for row 4 (copy formatting from Q4) ---
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 0.057927248158369672

# --- Update the selected/active cell as recorded in the saved view ---
$ws.Range("O10").Select()
